$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) - update F2, F5, F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7611
$ws1.Range("F5").Value = 25
$ws1.Range("F6").Value = 269

# Sheet "全部类型" (4th sheet) - update F2, F5, F6 (mirrors same data)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7611
$ws4.Range("F5").Value = 25
$ws4.Range("F6").Value = 269
